$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the example/placeholder bank-detail values (B18, B19, B20, B21, B28, B29).
# Clearing the only reference to these shared strings causes them (and the
# orphaned rich-text run in "Your Address...") to be dropped from the
# shared strings table when the workbook is rewritten.
$ws.Range("B18").Value = ""
$ws.Range("B19").Value = ""
$ws.Range("B20").Value = ""
$ws.Range("B21").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("B29").Value = ""

# Narrow the two big merged "value" blocks from column E to column D, and
# give row 28's B:D block a new centered style while doing it.
$ws.Range("B21:E27").UnMerge()
$ws.Range("B29:E35").UnMerge()

# Re-sequence the whole mergeCells list back to its natural top-to-bottom
# order (unmerging/merging appends to the end of the internal list).
$ws.Range("A1:E1").UnMerge()
$ws.Range("A5:E5").UnMerge()
$ws.Range("A6:E6").UnMerge()
$ws.Range("A7:E7").UnMerge()
$ws.Range("A10:E10").UnMerge()
$ws.Range("B15:E15").UnMerge()
$ws.Range("A16:E16").UnMerge()
$ws.Range("B18:E18").UnMerge()
$ws.Range("B19:E19").UnMerge()
$ws.Range("B20:E20").UnMerge()
$ws.Range("A21:A27").UnMerge()
$ws.Range("B28:D28").UnMerge()
$ws.Range("A29:A35").UnMerge()

$ws.Range("A1:E1").Merge()
$ws.Range("A5:E5").Merge()
$ws.Range("A6:E6").Merge()
$ws.Range("A7:E7").Merge()
$ws.Range("A10:E10").Merge()
$ws.Range("B15:E15").Merge()
$ws.Range("A16:E16").Merge()
$ws.Range("B18:E18").Merge()
$ws.Range("B19:E19").Merge()
$ws.Range("B20:E20").Merge()
$ws.Range("A21:A27").Merge()
$ws.Range("B21:D27").Merge()
$ws.Range("B28:D28").Merge()
$ws.Range("A29:A35").Merge()
$ws.Range("B29:D35").Merge()

# B28:D28 gets a new centered alignment style (a fresh cellXf distinct from
# the surrounding left/bottom one used by A28).
$ws.Range("B28:D28").HorizontalAlignment = -4108  ## xlCenter
$ws.Range("B28:D28").VerticalAlignment = -4108    ## xlCenter

# Selection moves to D9.
$ws.Range("D9").Select()
